$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11, columns B:G (ME, MAE, MSE, RMSE, SE, N)
$data = @(
    @(0.103921568627451,  0.7,                  0.9849019607843137, 0.9924222693915699, 0.9967869951627029, 51),
    @(0.01600000000000001,0.704,                0.9964,              0.9981983770774224, 1.008203089356383,  50),
    @(0.1081632653061224, 0.7285714285714285,   0.9793877551020409, 0.9896402149781712, 0.99390574978505,   49),
    @(0.006249999999999996,0.6687500000000001,  0.8072916666666666, 0.898494110535326,  0.9079802814430288, 48),
    @(0.0936170212765957, 0.7234042553191488,   0.9825531914893616, 0.9912382112738398, 0.9974760006525926, 47),
    @(0.004347826086956521,0.7304347826086957,  0.9669565217391302, 0.9833394743114557, 0.9941957155273536, 46),
    @(0.1,                0.7444444444444445,   0.918,               0.9581231653602788, 0.9636578042211683, 45),
    @(0.01136363636363638,0.7522727272727273,   1.020681818181818,  1.010287987744989,  1.021903356166623,  44),
    @(0.04186046511627906,0.7627906976744185,   0.9251162790697675, 0.9618296517937921, 0.9722905107061892, 43),
    @(0.05238095238095236,0.7,                  0.8742857142857141, 0.9350324669687754, 0.9448804532044524, 42)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("F$row").Value = $vals[4]
    $ws.Range("G$row").Value = $vals[5]
}

$wb.Save()
